$wb = $excel.ActiveWorkbook

# --- POBasedInvoice sheet: refresh invoice number / quantity / amount test data ---
$poSheet = $wb.Worksheets.Item("POBasedInvoice")

# Quantity (B), IGST (C) and Quantity-repeat (J) columns store numeric-looking
# values as text in this workbook, so force text formatting before assigning
# the literal values (keeps cells as shared strings, matching the test-data
# generator's original output).
$poSheet.Range("B2:C9").NumberFormat = "@"
$poSheet.Range("J2:J9").NumberFormat = "@"

$poSheet.Range("A2").Value = "TESTINV97192"
$poSheet.Range("B2").Value = "8"
$poSheet.Range("C2").Value = "1.44"
$poSheet.Range("J2").Value = "8"

$poSheet.Range("A3").Value = "TESTINV04156"
$poSheet.Range("B3").Value = "8"
$poSheet.Range("C3").Value = "1.44"
$poSheet.Range("J3").Value = "8"

$poSheet.Range("A4").Value = "TESTINV17714"
$poSheet.Range("B4").Value = "8"
$poSheet.Range("C4").Value = "1.44"
$poSheet.Range("J4").Value = "8"

$poSheet.Range("A5").Value = "TESTINV73341"
$poSheet.Range("B5").Value = "1"
$poSheet.Range("C5").Value = "0.18"
$poSheet.Range("J5").Value = "1"

$poSheet.Range("A6").Value = "TESTINV34048"
$poSheet.Range("B6").Value = "2"
$poSheet.Range("C6").Value = "0.36"
$poSheet.Range("J6").Value = "2"

$poSheet.Range("A7").Value = "TESTINV25913"
$poSheet.Range("B7").Value = "3"
$poSheet.Range("C7").Value = "0.54"
$poSheet.Range("J7").Value = "3"

$poSheet.Range("A8").Value = "TESTINV06532"
$poSheet.Range("B8").Value = "1"
$poSheet.Range("C8").Value = "0.18"
$poSheet.Range("J8").Value = "1"

$poSheet.Range("A9").Value = "TESTINV09872"
$poSheet.Range("B9").Value = "2"
$poSheet.Range("C9").Value = "0.36"
$poSheet.Range("J9").Value = "2"

$poSheet.Range("L14").Select()

# --- BADashboardPage sheet: add Submitting At / Submitting To columns ---
$dashSheet = $wb.Worksheets.Item("BADashboardPage")

$dashSheet.Range("C1").Value = "Submitting At"
$dashSheet.Range("D1").Value = "Submitting To"
$dashSheet.Range("C2").Value = "Ahmedabad"
$dashSheet.Range("D2").Value = "Nishant Gore"

$dashSheet.Columns.Item(3).ColumnWidth = 15.5546875
$dashSheet.Columns.Item(4).ColumnWidth = 14.44140625

$dashSheet.Range("D5").Select()
